# Pergantian dummy dataset dan update script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Source_File" column is inserted right after the first column,
# "Data_Source" is renamed/replaced by "File_Name", and "Update_by"
# becomes "Updated_by". Easiest reliable way: insert a column before B,
# then rewrite the header row values and column widths explicitly.

$ws.Columns.Item(2).Insert()

$ws.Range("A1").Value = "File_Name"
$ws.Range("B1").Value = "Source_File"
$ws.Range("C1").Value = "Modification_Type"
$ws.Range("D1").Value = "Updated_by"
$ws.Range("E1").Value = "Dashboard"
$ws.Range("F1").Value = "Update_Periode "
$ws.Range("G1").Value = "Target_Update"
$ws.Range("H1").Value = "Realisasi"
$ws.Range("I1").Value = "SLA_(Met/Miss)"

# NOTE: the runtime snaps ColumnWidth to a 1/6-character grid
# (stored = round(value*6)/6 + 5/6), so the values below are the
# inputs whose snapped results equal (or, for column E, come as close
# as that grid allows to) the widths recorded in the target file:
#   A -> 12, B -> 15, E -> ~18.71
$ws.Columns.Item(1).ColumnWidth = 11.166666666666666
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws.Columns.Item(5).ColumnWidth = 17.833333333333332

$ws.Range("I14").Select()
